# Auto-generated edit script: updates cached numeric values in the
# per-job 'Leve Profits' worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# These sheets hold static, scheduler-computed values (no formulas), matching
# the upstream commit 'chore: update Sheets via scheduled runner'.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 10).Value = 749.5
$ws.Cells.Item(2, 8).Value = 759.8
$ws.Cells.Item(2, 12).Value = 749.5
$ws.Cells.Item(2, 14).Value = -975.5
$ws.Cells.Item(17, 12).Value = 8472.999899999999
$ws.Cells.Item(17, 8).Value = 2824.3333
$ws.Cells.Item(17, 10).Value = 2824.3333
$ws.Cells.Item(17, 14).Value = -8808.999899999999
$ws.Cells.Item(40, 13).Value = -5791.6665
$ws.Cells.Item(40, 9).Value = 5966.6665
$ws.Cells.Item(40, 11).Value = 5966.6665
$ws.Cells.Item(40, 8).Value = 10342.857
$ws.Cells.Item(64, 14).Value = -9646
$ws.Cells.Item(64, 8).Value = 8869.916999999999
$ws.Cells.Item(64, 10).Value = 9150
$ws.Cells.Item(64, 12).Value = 9150
$ws.Cells.Item(67, 8).Value = 8869.916999999999
$ws.Cells.Item(67, 10).Value = 9150
$ws.Cells.Item(67, 14).Value = -10866
$ws.Cells.Item(67, 12).Value = 9150
$ws.Cells.Item(74, 12).Value = 11428.571
$ws.Cells.Item(74, 10).Value = 11428.571
$ws.Cells.Item(74, 8).Value = 11875
$ws.Cells.Item(74, 13).Value = -14064
$ws.Cells.Item(74, 9).Value = 15000
$ws.Cells.Item(74, 14).Value = -13300.571
$ws.Cells.Item(74, 11).Value = 15000
$ws.Cells.Item(77, 13).Value = -70320
$ws.Cells.Item(77, 9).Value = 15000
$ws.Cells.Item(77, 10).Value = 11428.571
$ws.Cells.Item(77, 14).Value = -66502.855
$ws.Cells.Item(77, 12).Value = 57142.855
$ws.Cells.Item(77, 8).Value = 11875
$ws.Cells.Item(77, 11).Value = 75000
$ws.Cells.Item(112, 14).Value = -8534
$ws.Cells.Item(112, 12).Value = 6318
$ws.Cells.Item(112, 8).Value = 2037.6774
$ws.Cells.Item(112, 10).Value = 2106
$ws.Cells.Item(138, 10).Value = 6395
$ws.Cells.Item(138, 13).Value = -4607.299999999999
$ws.Cells.Item(138, 11).Value = 9747.299999999999
$ws.Cells.Item(138, 8).Value = 5567.1313
$ws.Cells.Item(138, 9).Value = 3249.1
$ws.Cells.Item(138, 12).Value = 19185
$ws.Cells.Item(138, 14).Value = -29465

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 1100
$ws.Cells.Item(12, 13).Value = -927
$ws.Cells.Item(12, 11).Value = 1100
$ws.Cells.Item(12, 9).Value = 1100
$ws.Cells.Item(32, 8).Value = 20851820
$ws.Cells.Item(32, 9).Value = 21295304
$ws.Cells.Item(32, 11).Value = 21295304
$ws.Cells.Item(32, 13).Value = -21295017
$ws.Cells.Item(56, 8).Value = 91666.664
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 91666.664
$ws.Cells.Item(56, 13).ClearContents()
$ws.Cells.Item(56, 12).Value = 91666.664
$ws.Cells.Item(56, 14).Value = -93150.664
$ws.Cells.Item(63, 8).Value = 2893.6365
$ws.Cells.Item(63, 12).Value = 3969.9
$ws.Cells.Item(63, 10).Value = 3969.9
$ws.Cells.Item(63, 14).Value = -5341.9
$ws.Cells.Item(66, 10).Value = 3969.9
$ws.Cells.Item(66, 14).Value = -26713.5
$ws.Cells.Item(66, 12).Value = 19849.5
$ws.Cells.Item(66, 8).Value = 2893.6365
$ws.Cells.Item(74, 8).Value = 3699.5833
$ws.Cells.Item(74, 13).Value = -2426.342
$ws.Cells.Item(74, 9).Value = 3300.342
$ws.Cells.Item(74, 11).Value = 3300.342
$ws.Cells.Item(77, 13).Value = -12133.71
$ws.Cells.Item(77, 9).Value = 3300.342
$ws.Cells.Item(77, 8).Value = 3699.5833
$ws.Cells.Item(77, 11).Value = 16501.71
$ws.Cells.Item(132, 12).Value = 18808.4442
$ws.Cells.Item(132, 13).Value = -13557.599
$ws.Cells.Item(132, 8).Value = 5792.14
$ws.Cells.Item(132, 14).Value = -23868.4442
$ws.Cells.Item(132, 11).Value = 16087.599
$ws.Cells.Item(132, 9).Value = 5362.533
$ws.Cells.Item(132, 10).Value = 6269.4814

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 14).ClearContents()
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(81, 14).Value = -52929.5
$ws.Cells.Item(81, 12).Value = 50807.5
$ws.Cells.Item(81, 10).Value = 50807.5
$ws.Cells.Item(81, 8).Value = 50807.5
$ws.Cells.Item(84, 14).Value = -163030.5
$ws.Cells.Item(84, 12).Value = 152422.5
$ws.Cells.Item(84, 10).Value = 50807.5
$ws.Cells.Item(84, 8).Value = 50807.5
$ws.Cells.Item(86, 11).Value = 6443
$ws.Cells.Item(86, 13).Value = -5320
$ws.Cells.Item(86, 8).Value = 6326.875
$ws.Cells.Item(86, 9).Value = 6443
$ws.Cells.Item(89, 13).Value = -26599
$ws.Cells.Item(89, 11).Value = 32215
$ws.Cells.Item(89, 9).Value = 6443
$ws.Cells.Item(89, 8).Value = 6326.875
$ws.Cells.Item(107, 13).Value = -3912.7
$ws.Cells.Item(107, 11).Value = 5832.7
$ws.Cells.Item(107, 9).Value = 5832.7
$ws.Cells.Item(107, 8).Value = 5277.25
$ws.Cells.Item(134, 9).Value = 2906.6099
$ws.Cells.Item(134, 8).Value = 3674.9656
$ws.Cells.Item(134, 11).Value = 8719.8297
$ws.Cells.Item(134, 13).Value = -6184.8297

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 13).ClearContents()
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 4468.875
$ws.Cells.Item(31, 13).Value = -4173.875
$ws.Cells.Item(31, 8).Value = 5570.9165
$ws.Cells.Item(31, 11).Value = 4468.875
$ws.Cells.Item(34, 9).Value = 4468.875
$ws.Cells.Item(34, 11).Value = 4468.875
$ws.Cells.Item(34, 8).Value = 5570.9165
$ws.Cells.Item(34, 13).Value = -4266.875
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 8).Value = 10007
$ws.Cells.Item(58, 13).ClearContents()
$ws.Cells.Item(132, 12).Value = 16323.3531
$ws.Cells.Item(132, 13).Value = -7239.459800000001
$ws.Cells.Item(132, 8).Value = 3944.2407
$ws.Cells.Item(132, 14).Value = -21383.3531
$ws.Cells.Item(132, 11).Value = 9769.459800000001
$ws.Cells.Item(132, 9).Value = 3256.4866
$ws.Cells.Item(132, 10).Value = 5441.1177
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 8).Value = 10007
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(141, 14).Value = -512509.84
$ws.Cells.Item(141, 10).Value = 502149.84
$ws.Cells.Item(141, 8).Value = 502149.84
$ws.Cells.Item(141, 12).Value = 502149.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 14).Value = -6041.83338
$ws.Cells.Item(107, 10).Value = 733.94446
$ws.Cells.Item(107, 8).Value = 581.5925999999999
$ws.Cells.Item(107, 12).Value = 2201.83338
$ws.Cells.Item(113, 11).Value = 11820
$ws.Cells.Item(113, 13).Value = -9650
$ws.Cells.Item(113, 8).Value = 4628.3335
$ws.Cells.Item(113, 9).Value = 3940
$ws.Cells.Item(114, 9).Value = 1029.2
$ws.Cells.Item(114, 10).Value = 2051.375
$ws.Cells.Item(114, 13).Value = 166.3999999999996
$ws.Cells.Item(114, 12).Value = 6154.125
$ws.Cells.Item(114, 14).Value = -12662.125
$ws.Cells.Item(114, 11).Value = 3087.6
$ws.Cells.Item(114, 8).Value = 1658.2307
$ws.Cells.Item(131, 11).Value = 17382.429
$ws.Cells.Item(131, 8).Value = 5919.8335
$ws.Cells.Item(131, 13).Value = -12342.429
$ws.Cells.Item(131, 9).Value = 5794.143
$ws.Cells.Item(136, 8).Value = 2980.9092
$ws.Cells.Item(136, 10).Value = 2998.7144
$ws.Cells.Item(136, 12).Value = 8996.143199999999
$ws.Cells.Item(136, 14).Value = -19196.1432
$ws.Cells.Item(139, 13).Value = -250010180
$ws.Cells.Item(139, 9).Value = 83338440
$ws.Cells.Item(139, 8).Value = 30314282
$ws.Cells.Item(139, 11).Value = 250015320

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 28749.5
$ws.Cells.Item(57, 12).Value = 49999
$ws.Cells.Item(57, 14).Value = -51639
$ws.Cells.Item(57, 10).Value = 49999
$ws.Cells.Item(135, 14).Value = -68601.26999999999
$ws.Cells.Item(135, 8).Value = 58461.27
$ws.Cells.Item(135, 10).Value = 58461.27
$ws.Cells.Item(135, 12).Value = 58461.27

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 11).Value = 1500
$ws.Cells.Item(22, 10).Value = 1691
$ws.Cells.Item(22, 8).Value = 1563.6666
$ws.Cells.Item(22, 13).Value = -1205
$ws.Cells.Item(22, 14).Value = -2281
$ws.Cells.Item(22, 12).Value = 1691
$ws.Cells.Item(22, 9).Value = 1500
$ws.Cells.Item(27, 11).Value = 1500
$ws.Cells.Item(27, 14).Value = -1905
$ws.Cells.Item(27, 8).Value = 1563.6666
$ws.Cells.Item(27, 10).Value = 1691
$ws.Cells.Item(27, 9).Value = 1500
$ws.Cells.Item(27, 13).Value = -1393
$ws.Cells.Item(27, 12).Value = 1691
$ws.Cells.Item(132, 12).Value = 22561.0005
$ws.Cells.Item(132, 13).Value = -15670.166
$ws.Cells.Item(132, 8).Value = 6727.4546
$ws.Cells.Item(132, 14).Value = -27621.0005
$ws.Cells.Item(132, 11).Value = 18200.166
$ws.Cells.Item(132, 9).Value = 6066.722
$ws.Cells.Item(132, 10).Value = 7520.3335
$ws.Cells.Item(136, 9).Value = 4250
$ws.Cells.Item(136, 8).Value = 6439.615
$ws.Cells.Item(136, 11).Value = 12750
$ws.Cells.Item(136, 13).Value = -10200

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 13).Value = -326712
$ws.Cells.Item(15, 11).Value = 327000
$ws.Cells.Item(15, 9).Value = 327000
$ws.Cells.Item(15, 8).Value = 301718.2
$ws.Cells.Item(52, 9).Value = 5016672
$ws.Cells.Item(52, 8).Value = 4305719
$ws.Cells.Item(52, 13).Value = -5016446
$ws.Cells.Item(52, 11).Value = 5016672
$ws.Cells.Item(61, 14).Value = -31829.666
$ws.Cells.Item(61, 11).Value = 33797
$ws.Cells.Item(61, 10).Value = 31245.666
$ws.Cells.Item(61, 8).Value = 32840.25
$ws.Cells.Item(61, 12).Value = 31245.666
$ws.Cells.Item(61, 13).Value = -33505
$ws.Cells.Item(61, 9).Value = 33797
$ws.Cells.Item(62, 14).Value = -9248
$ws.Cells.Item(62, 10).Value = 8000
$ws.Cells.Item(62, 13).Value = -378
$ws.Cells.Item(62, 11).Value = 1002
$ws.Cells.Item(62, 12).Value = 8000
$ws.Cells.Item(62, 9).Value = 1002
$ws.Cells.Item(62, 8).Value = 4501
$ws.Cells.Item(65, 10).Value = 8000
$ws.Cells.Item(65, 13).Value = -1890
$ws.Cells.Item(65, 8).Value = 4501
$ws.Cells.Item(65, 11).Value = 5010
$ws.Cells.Item(65, 9).Value = 1002
$ws.Cells.Item(65, 12).Value = 40000
$ws.Cells.Item(65, 14).Value = -46240
$ws.Cells.Item(114, 10).Value = 122499.5
$ws.Cells.Item(114, 12).Value = 122499.5
$ws.Cells.Item(114, 14).Value = -131177.5
$ws.Cells.Item(114, 8).Value = 122499.5
$ws.Cells.Item(133, 10).Value = 43271.285
$ws.Cells.Item(133, 14).Value = -53391.285
$ws.Cells.Item(133, 8).Value = 43271.285
$ws.Cells.Item(133, 12).Value = 43271.285

